$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value = 0
$ws.Cells.Item(64, 9).Value = 0
$ws.Cells.Item(64, 10).Value = 0
$ws.Cells.Item(64, 11).Value = 0
$ws.Cells.Item(64, 12).Value = 0
$ws.Cells.Item(64, 14).Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(67, 8).Value = 0
$ws.Cells.Item(67, 9).Value = 0
$ws.Cells.Item(67, 10).Value = 0
$ws.Cells.Item(67, 11).Value = 0
$ws.Cells.Item(67, 12).Value = 0
$ws.Cells.Item(67, 14).Value = ""

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 11299.6
$ws.Cells.Item(132, 9).Value = 10499
$ws.Cells.Item(132, 10).Value = 11833.333
$ws.Cells.Item(132, 11).Value = 31497
$ws.Cells.Item(132, 12).Value = 35499.999
$ws.Cells.Item(132, 13).Value = -28967

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(137, 8).Value = 3981.0417
$ws.Cells.Item(137, 9).Value = 3999.3333
$ws.Cells.Item(137, 10).Value = 3974.9443
$ws.Cells.Item(137, 11).Value = 11997.9999
$ws.Cells.Item(137, 12).Value = 11924.8329
$ws.Cells.Item(137, 13).Value = -9447.999899999999
$ws.Cells.Item(137, 14).Value = -17024.8329

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 6468.8887
$ws.Cells.Item(138, 9).Value = 7344.2
$ws.Cells.Item(138, 10).Value = 5374.75
$ws.Cells.Item(138, 11).Value = 22032.6
$ws.Cells.Item(138, 12).Value = 16124.25
$ws.Cells.Item(138, 13).Value = -16892.6
$ws.Cells.Item(138, 14).Value = -26404.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 300.5
$ws.Cells.Item(4, 9).Value = 300.5
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 300.5
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = -184.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14063.143
$ws.Cells.Item(32, 9).Value = 11690.4
$ws.Cells.Item(32, 10).Value = 19995
$ws.Cells.Item(32, 11).Value = 11690.4
$ws.Cells.Item(32, 12).Value = 19995
$ws.Cells.Item(32, 13).Value = -11403.4
$ws.Cells.Item(32, 14).Value = -20569

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 3249.5
$ws.Cells.Item(61, 9).Value = 3249.5
$ws.Cells.Item(61, 10).Value = 0
$ws.Cells.Item(61, 11).Value = 3249.5
$ws.Cells.Item(61, 12).Value = 0
$ws.Cells.Item(61, 13).Value = -3037.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 3272.25
$ws.Cells.Item(63, 9).Value = 3272.25
$ws.Cells.Item(63, 10).Value = 0
$ws.Cells.Item(63, 11).Value = 3272.25
$ws.Cells.Item(63, 12).Value = 0
$ws.Cells.Item(63, 13).Value = -2586.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(66, 8).Value = 3272.25
$ws.Cells.Item(66, 9).Value = 3272.25
$ws.Cells.Item(66, 10).Value = 0
$ws.Cells.Item(66, 11).Value = 16361.25
$ws.Cells.Item(66, 12).Value = 0
$ws.Cells.Item(66, 13).Value = -12929.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(88, 8).Value = 1100
$ws.Cells.Item(88, 9).Value = 1100
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 11).Value = 1100
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(88, 13).Value = -694
$ws.Cells.Item(88, 14).Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(91, 8).Value = 1100
$ws.Cells.Item(91, 9).Value = 1100
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 11).Value = 1100
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(91, 13).Value = 304
$ws.Cells.Item(91, 14).Value = ""

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(132, 8).Value = 6001.1
$ws.Cells.Item(132, 9).Value = 6001.1
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 18003.3
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -15473.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(136, 8).Value = 3249.5
$ws.Cells.Item(136, 9).Value = 3249.5
$ws.Cells.Item(136, 10).Value = 0
$ws.Cells.Item(136, 11).Value = 9748.5
$ws.Cells.Item(136, 12).Value = 0
$ws.Cells.Item(136, 13).Value = -7198.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(103, 8).Value = 8575
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 10).Value = 8575
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 12).Value = 8575
$ws.Cells.Item(103, 14).Value = -10919

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 2268.1667
$ws.Cells.Item(105, 9).Value = 2402
$ws.Cells.Item(105, 10).Value = 1599
$ws.Cells.Item(105, 11).Value = 2402
$ws.Cells.Item(105, 12).Value = 1599
$ws.Cells.Item(105, 13).Value = -655

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(107, 8).Value = 11093.833
$ws.Cells.Item(107, 9).Value = 11093.833
$ws.Cells.Item(107, 10).Value = 0
$ws.Cells.Item(107, 11).Value = 11093.833
$ws.Cells.Item(107, 12).Value = 0
$ws.Cells.Item(107, 13).Value = -9173.833000000001
$ws.Cells.Item(107, 14).Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 2681.5
$ws.Cells.Item(99, 9).Value = 2681.5
$ws.Cells.Item(99, 10).Value = 0
$ws.Cells.Item(99, 11).Value = 2681.5
$ws.Cells.Item(99, 12).Value = 0
$ws.Cells.Item(99, 13).Value = -1183.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(126, 8).Value = 2681.5
$ws.Cells.Item(126, 9).Value = 2681.5
$ws.Cells.Item(126, 10).Value = 0
$ws.Cells.Item(126, 11).Value = 8044.5
$ws.Cells.Item(126, 12).Value = 0
$ws.Cells.Item(126, 13).Value = -5574.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(132, 8).Value = 3200.9
$ws.Cells.Item(132, 9).Value = 3131.125
$ws.Cells.Item(132, 10).Value = 3480
$ws.Cells.Item(132, 11).Value = 9393.375
$ws.Cells.Item(132, 12).Value = 10440
$ws.Cells.Item(132, 13).Value = -6863.375
$ws.Cells.Item(132, 14).Value = -15500

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 125253
$ws.Cells.Item(4, 9).Value = 325
$ws.Cells.Item(4, 10).Value = 333466.34
$ws.Cells.Item(4, 11).Value = 975
$ws.Cells.Item(4, 12).Value = 1000399.02
$ws.Cells.Item(4, 13).Value = -863
$ws.Cells.Item(4, 14).Value = -1000623.02

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 156.11111
$ws.Cells.Item(7, 9).Value = 113.125
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 339.375
$ws.Cells.Item(7, 12).Value = 1500
$ws.Cells.Item(7, 13).Value = -227.375
$ws.Cells.Item(7, 14).Value = -1724

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(15, 8).Value = 80
$ws.Cells.Item(15, 9).Value = 80
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 240
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = -100

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(40, 8).Value = 305.30768
$ws.Cells.Item(40, 9).Value = 99
$ws.Cells.Item(40, 10).Value = 342.81818
$ws.Cells.Item(40, 11).Value = 396
$ws.Cells.Item(40, 12).Value = 1371.27272
$ws.Cells.Item(40, 13).Value = -327
$ws.Cells.Item(40, 14).Value = -1509.27272

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(21, 8).Value = 1439214.2
$ws.Cells.Item(21, 9).Value = 5000000
$ws.Cells.Item(21, 10).Value = 14900
$ws.Cells.Item(21, 11).Value = 5000000
$ws.Cells.Item(21, 12).Value = 14900
$ws.Cells.Item(21, 13).Value = -4999827
$ws.Cells.Item(21, 14).Value = -15246

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(30, 8).Value = 1439214.2
$ws.Cells.Item(30, 9).Value = 5000000
$ws.Cells.Item(30, 10).Value = 14900
$ws.Cells.Item(30, 11).Value = 5000000
$ws.Cells.Item(30, 12).Value = 14900
$ws.Cells.Item(30, 13).Value = -4999895
$ws.Cells.Item(30, 14).Value = -15110

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(62, 8).Value = 0
$ws.Cells.Item(62, 9).Value = 0
$ws.Cells.Item(62, 10).Value = 0
$ws.Cells.Item(62, 11).Value = 0
$ws.Cells.Item(62, 12).Value = 0
$ws.Cells.Item(62, 14).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(65, 8).Value = 0
$ws.Cells.Item(65, 9).Value = 0
$ws.Cells.Item(65, 10).Value = 0
$ws.Cells.Item(65, 11).Value = 0
$ws.Cells.Item(65, 12).Value = 0
$ws.Cells.Item(65, 14).Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 5198.8
$ws.Cells.Item(102, 9).Value = 5198.8
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 11).Value = 5198.8
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 13).Value = -3576.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 8117.9443
$ws.Cells.Item(122, 9).Value = 3017.0625
$ws.Cells.Item(122, 10).Value = 48925
$ws.Cells.Item(122, 11).Value = 9051.1875
$ws.Cells.Item(122, 12).Value = 146775
$ws.Cells.Item(122, 13).Value = -6601.1875

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 2033.5
$ws.Cells.Item(132, 9).Value = 2033.5
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 11).Value = 6100.5
$ws.Cells.Item(132, 12).Value = 0
$ws.Cells.Item(132, 13).Value = -3570.5
$ws.Cells.Item(132, 14).Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7336.625
$ws.Cells.Item(40, 9).Value = 7115.5
$ws.Cells.Item(40, 10).Value = 8000
$ws.Cells.Item(40, 11).Value = 7115.5
$ws.Cells.Item(40, 12).Value = 8000
$ws.Cells.Item(40, 13).Value = -6979.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(41, 8).Value = 5000
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 10).Value = 5000
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 12).Value = 5000
$ws.Cells.Item(41, 13).Value = ""
$ws.Cells.Item(41, 14).Value = -5876

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(42, 8).Value = 35640
$ws.Cells.Item(42, 9).Value = 30025
$ws.Cells.Item(42, 10).Value = 41255
$ws.Cells.Item(42, 11).Value = 30025
$ws.Cells.Item(42, 12).Value = 41255
$ws.Cells.Item(42, 13).Value = -29462
$ws.Cells.Item(42, 14).Value = -42381

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(49, 8).Value = 35640
$ws.Cells.Item(49, 9).Value = 30025
$ws.Cells.Item(49, 10).Value = 41255
$ws.Cells.Item(49, 11).Value = 30025
$ws.Cells.Item(49, 12).Value = 41255
$ws.Cells.Item(49, 13).Value = -29878
$ws.Cells.Item(49, 14).Value = -41549

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 2674
$ws.Cells.Item(93, 9).Value = 2770.2856
$ws.Cells.Item(93, 10).Value = 2000
$ws.Cells.Item(93, 11).Value = 2770.2856
$ws.Cells.Item(93, 12).Value = 2000
$ws.Cells.Item(93, 13).Value = -1522.2856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(100, 8).Value = 4488.8887
$ws.Cells.Item(100, 9).Value = 3900
$ws.Cells.Item(100, 10).Value = 5666.6665
$ws.Cells.Item(100, 11).Value = 3900
$ws.Cells.Item(100, 12).Value = 5666.6665
$ws.Cells.Item(100, 13).Value = -3359
$ws.Cells.Item(100, 14).Value = -6748.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(132, 8).Value = 21866.182
$ws.Cells.Item(132, 9).Value = 19052.8
$ws.Cells.Item(132, 10).Value = 50000
$ws.Cells.Item(132, 11).Value = 57158.39999999999
$ws.Cells.Item(132, 12).Value = 150000
$ws.Cells.Item(132, 13).Value = -54628.39999999999
$ws.Cells.Item(132, 14).Value = -155060

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 3750
$ws.Cells.Item(62, 9).Value = 3500
$ws.Cells.Item(62, 10).Value = 4000
$ws.Cells.Item(62, 11).Value = 3500
$ws.Cells.Item(62, 12).Value = 4000
$ws.Cells.Item(62, 13).Value = -2876
$ws.Cells.Item(62, 14).Value = -5248

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(65, 8).Value = 3750
$ws.Cells.Item(65, 9).Value = 3500
$ws.Cells.Item(65, 10).Value = 4000
$ws.Cells.Item(65, 11).Value = 17500
$ws.Cells.Item(65, 12).Value = 20000
$ws.Cells.Item(65, 13).Value = -14380
$ws.Cells.Item(65, 14).Value = -26240

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 0
$ws.Cells.Item(68, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 0
$ws.Cells.Item(71, 14).Value = ""

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(122, 8).Value = 501596.5
$ws.Cells.Item(122, 9).Value = 667462.3
$ws.Cells.Item(122, 10).Value = 3999
$ws.Cells.Item(122, 11).Value = 2002386.9
$ws.Cells.Item(122, 12).Value = 11997
$ws.Cells.Item(122, 13).Value = -1999936.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 3533.1
$ws.Cells.Item(136, 9).Value = 3166.5
$ws.Cells.Item(136, 10).Value = 4999.5
$ws.Cells.Item(136, 11).Value = 9499.5
$ws.Cells.Item(136, 12).Value = 14998.5
$ws.Cells.Item(136, 13).Value = -6949.5
$ws.Cells.Item(136, 14).Value = -20098.5
